$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "dnasr281@gmail.com, System"
$replacement = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
